$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3564.3704
$ws.Range("J40").Value = 3995
$ws.Range("L40").Value = 3995
$ws.Range("N40").Value = -4345
$ws.Range("H58").Value = 1121.1428
$ws.Range("I58").Value = 966
$ws.Range("J58").Value = 1237.5
$ws.Range("K58").Value = 2898
$ws.Range("L58").Value = 3712.5
$ws.Range("M58").Value = -2748
$ws.Range("N58").Value = -4012.5
$ws.Range("H86").Value = 14514.4
$ws.Range("I86").Value = 14793.125
$ws.Range("J86").Value = 13399.5
$ws.Range("K86").Value = 14793.125
$ws.Range("L86").Value = 13399.5
$ws.Range("M86").Value = -13670.125
$ws.Range("N86").Value = -15645.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H89").Value = 14514.4
$ws.Range("I89").Value = 14793.125
$ws.Range("J89").Value = 13399.5
$ws.Range("K89").Value = 73965.625
$ws.Range("L89").Value = 66997.5
$ws.Range("M89").Value = -68349.625
$ws.Range("N89").Value = -78229.5
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H127").Value = 1615.2727
$ws.Range("I127").Value = 1526.8
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 4580.4
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = 379.6000000000004
$ws.Range("N127").Value = -17420
$ws.Range("H132").Value = 1486.1428
$ws.Range("I132").Value = 1477.069
$ws.Range("K132").Value = 4431.207
$ws.Range("M132").Value = -1901.207
$ws.Range("H137").Value = 3259.1462
$ws.Range("I137").Value = 1660.5186
$ws.Range("K137").Value = 4981.5558
$ws.Range("M137").Value = -2431.5558
$ws.Range("H138").Value = 2050.8794
$ws.Range("I138").Value = 881.05884
$ws.Range("J138").Value = 2535.9268
$ws.Range("K138").Value = 2643.17652
$ws.Range("L138").Value = 7607.780400000001
$ws.Range("M138").Value = 2496.82348
$ws.Range("N138").Value = -17887.7804

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2773.9092
$ws.Range("I45").Value = 1505
$ws.Range("J45").Value = 3249.75
$ws.Range("K45").Value = 1505
$ws.Range("L45").Value = 3249.75
$ws.Range("M45").Value = -1128
$ws.Range("N45").Value = -4003.75
$ws.Range("H61").Value = 14742685
$ws.Range("I61").Value = 18524244
$ws.Range("K61").Value = 18524244
$ws.Range("M61").Value = -18524032
$ws.Range("H63").Value = 3443.2
$ws.Range("I63").Value = 2765.75
$ws.Range("K63").Value = 2765.75
$ws.Range("M63").Value = -2079.75
$ws.Range("H66").Value = 3443.2
$ws.Range("I66").Value = 2765.75
$ws.Range("K66").Value = 13828.75
$ws.Range("M66").Value = -10396.75
$ws.Range("H74").Value = 31274932
$ws.Range("I74").Value = 250000000
$ws.Range("J74").Value = 28493.143
$ws.Range("K74").Value = 250000000
$ws.Range("L74").Value = 28493.143
$ws.Range("M74").Value = -249999126
$ws.Range("N74").Value = -30241.143
$ws.Range("H77").Value = 31274932
$ws.Range("I77").Value = 250000000
$ws.Range("J77").Value = 28493.143
$ws.Range("K77").Value = 1250000000
$ws.Range("L77").Value = 142465.715
$ws.Range("M77").Value = -1249995632
$ws.Range("N77").Value = -151201.715
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H132").Value = 6939.6553
$ws.Range("I132").Value = 2907.95
$ws.Range("K132").Value = 8723.849999999999
$ws.Range("M132").Value = -6193.849999999999
$ws.Range("H136").Value = 14742685
$ws.Range("I136").Value = 18524244
$ws.Range("K136").Value = 55572732
$ws.Range("M136").Value = -55570182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 155135.42
$ws.Range("I134").Value = 10420
$ws.Range("J134").Value = 179254.67
$ws.Range("K134").Value = 31260
$ws.Range("L134").Value = 537764.01
$ws.Range("M134").Value = -28725
$ws.Range("N134").Value = -542834.01

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("H22").Value = 678.5714
$ws.Range("I22").Value = 678.5714
$ws.Range("K22").Value = 678.5714
$ws.Range("M22").Value = -328.5714
$ws.Range("H87").Value = 62109.832
$ws.Range("J87").Value = 88829.5
$ws.Range("L87").Value = 88829.5
$ws.Range("N87").Value = -91201.5
$ws.Range("H90").Value = 62109.832
$ws.Range("J90").Value = 88829.5
$ws.Range("L90").Value = 266488.5
$ws.Range("N90").Value = -278344.5
$ws.Range("H99").Value = 2570.8
$ws.Range("I99").Value = 2348.625
$ws.Range("J99").Value = 2824.7144
$ws.Range("K99").Value = 2348.625
$ws.Range("L99").Value = 2824.7144
$ws.Range("M99").Value = -850.625
$ws.Range("N99").Value = -5820.7144
$ws.Range("H105").Value = 2559.1365
$ws.Range("I105").Value = 1463.909
$ws.Range("K105").Value = 1463.909
$ws.Range("M105").Value = 283.0909999999999
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H126").Value = 2570.8
$ws.Range("I126").Value = 2348.625
$ws.Range("J126").Value = 2824.7144
$ws.Range("K126").Value = 7045.875
$ws.Range("L126").Value = 8474.143199999999
$ws.Range("M126").Value = -4575.875
$ws.Range("N126").Value = -13414.1432
$ws.Range("H134").Value = 2016007.8
$ws.Range("I134").Value = 10000000
$ws.Range("J134").Value = 20009.75
$ws.Range("K134").Value = 30000000
$ws.Range("L134").Value = 60029.25
$ws.Range("M134").Value = -29997465
$ws.Range("N134").Value = -65099.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2965.6667
$ws.Range("I5").Value = 400
$ws.Range("J5").Value = 4248.5
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 12745.5
$ws.Range("M5").Value = -1088
$ws.Range("N5").Value = -12969.5
$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224
$ws.Range("H113").Value = 2330
$ws.Range("I113").Value = 790
$ws.Range("J113").Value = 3100
$ws.Range("K113").Value = 2370
$ws.Range("L113").Value = 9300
$ws.Range("M113").Value = -200
$ws.Range("N113").Value = -13640
$ws.Range("H117").Value = 1390.2727
$ws.Range("J117").Value = 1658.8
$ws.Range("L117").Value = 4976.4
$ws.Range("N117").Value = -11860.4
$ws.Range("H122").Value = 2065.1428
$ws.Range("I122").Value = 250
$ws.Range("J122").Value = 2367.6667
$ws.Range("K122").Value = 2250
$ws.Range("L122").Value = 21309.0003
$ws.Range("M122").Value = 200
$ws.Range("N122").Value = -26209.0003
$ws.Range("H131").Value = 6200.6875
$ws.Range("I131").Value = 5621.8
$ws.Range("J131").Value = 7165.5
$ws.Range("K131").Value = 16865.4
$ws.Range("L131").Value = 21496.5
$ws.Range("M131").Value = -11825.4
$ws.Range("N131").Value = -31576.5
$ws.Range("H135").Value = 2965.6667
$ws.Range("I135").Value = 400
$ws.Range("J135").Value = 4248.5
$ws.Range("K135").Value = 3600
$ws.Range("L135").Value = 38236.5
$ws.Range("M135").Value = -1065
$ws.Range("N135").Value = -43306.5
$ws.Range("H140").Value = 217285.64
$ws.Range("I140").Value = 217285.64
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 651856.92
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -646676.92
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
$ws.Range("H102").Value = 5332.091
$ws.Range("I102").Value = 3294.7778
$ws.Range("K102").Value = 3294.7778
$ws.Range("M102").Value = -1672.7778
$ws.Range("H132").Value = 62517870
$ws.Range("I132").Value = 76929304
$ws.Range("K132").Value = 230787912
$ws.Range("M132").Value = -230785382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3402.2666
$ws.Range("I40").Value = 2759.5217
$ws.Range("K40").Value = 2759.5217
$ws.Range("M40").Value = -2623.5217
$ws.Range("H46").Value = 6360.1816
$ws.Range("J46").Value = 9916.666999999999
$ws.Range("L46").Value = 9916.666999999999
$ws.Range("N46").Value = -10292.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9708.666999999999
$ws.Range("J45").Value = 9708.666999999999
$ws.Range("L45").Value = 9708.666999999999
$ws.Range("N45").Value = -10690.667
$ws.Range("H108").Value = 113000
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").Value = ""
$ws.Range("H123").Value = 64410.832
$ws.Range("J123").Value = 64410.832
$ws.Range("L123").Value = 64410.832
$ws.Range("N123").Value = -74210.83199999999

